$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New inventory row (row 40): "889781" Cuchilla de limpieza (Kyocera KM 1500 1815 1820, FS 1000 1010 1018 1020 1050)
# Column A holds a numeric-looking product code that must stay TEXT (like every other
# code in column A, e.g. "9I7AQO"). A plain Range.Value assignment of "889781" would be
# auto-coerced to a number, so it is entered as a text-returning formula and then
# frozen to a static value via Copy/PasteSpecial(values) — this avoids leaving a
# formula behind AND avoids leaving a stray NumberFormat/quote-prefix style on the cell.
$ws.Range("A40").Formula = '="889781"'
$ws.Range("A40").Copy()
$ws.Range("A40").PasteSpecial(-4163)

$ws.Range("B40").Value = "Cuchilla de limpieza"
$ws.Range("C40").Value = "Kyocera KM 1500 1815 1820, FS 1000 1010 1018 1020 1050"
$ws.Range("D40").Value = 0
$ws.Range("E40").Value = 100000
$ws.Range("F40").Value = 1
$ws.Range("G40").Value = 0
$ws.Range("H40").Formula = "=(E40-D40)*G40"
$ws.Range("I40").Formula = "=D40*F40"
$ws.Range("J40").Value = 0
